$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns keep their original text representation
# (avoid Excel auto-converting numeric-looking strings to numbers and
# dropping significant trailing zeros / digit grouping).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '27.045.66'
$ws.Range("E2").Value = '  -2.82%  '
$ws.Range("D3").Value = '1.797.45'
$ws.Range("E3").Value = '  -3.16%  '
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").Value = '307.39'
$ws.Range("E5").Value = '  -2.97%  '
$ws.Range("E6").Value = '  -0.09%  '
$ws.Range("D7").Value = '0.4197'
$ws.Range("E7").Value = '  -2.96%  '
$ws.Range("E8").Value = '  -3.53%  '
$ws.Range("D9").Value = '0.07103'
$ws.Range("E9").Value = '  -3.93%  '
$ws.Range("D10").Value = '0.8445'
$ws.Range("E10").Value = '  -4.30%  '
$ws.Range("D11").Value = '20.13'
$ws.Range("E11").Value = '  -5.12%  '
$ws.Range("D12").Value = '1.805.18'
$ws.Range("E12").Value = '  -4.02%  '
$ws.Range("D13").Value = '5.295'
$ws.Range("E13").Value = '  -4.06%  '
$ws.Range("E14").Value = '  -4.18%  '
$ws.Range("D15").Value = '0.06755'
$ws.Range("E15").Value = '  -3.05%  '
$ws.Range("D16").Value = '1.005'
$ws.Range("E16").Value = '  +0.31%  '
$ws.Range("D17").Value = '80.32'
$ws.Range("E17").Value = '  -2.03%  '
$ws.Range("D18").Value = '0.000008697'
$ws.Range("E18").Value = '  -4.62%  '
$ws.Range("E19").Value = '  -0.04%  '
$ws.Range("D20").Value = '15.00'
$ws.Range("E20").Value = '  -4.15%  '
$ws.Range("D21").Value = '27.032.68'
$ws.Range("E21").Value = '  -3.11%  '
$ws.Range("D22").Value = '5.048'
$ws.Range("E22").Value = '  -1.06%  '
$ws.Range("D23").Value = '11.00'
$ws.Range("E23").Value = '  -0.40%  '
$ws.Range("D24").Value = '2.013.15'
$ws.Range("E24").Value = '  -4.06%  '
$ws.Range("D25").Value = '1.923'
$ws.Range("E25").Value = '  -3.12%  '
$ws.Range("D26").Value = '152.78'
$ws.Range("E26").Value = '  -1.36%  '
$ws.Range("D27").Value = '18.09'
$ws.Range("E27").Value = '  -5.72%  '
$ws.Range("D28").Value = '5.007'
$ws.Range("E28").Value = '  -6.80%  '
$ws.Range("D29").Value = '112.89'
$ws.Range("E29").Value = '  -3.04%  '
$ws.Range("E30").Value = '  -12.91%  '
$ws.Range("D31").Value = '0.09012'
$ws.Range("E31").Value = '  +0.51%  '
$ws.Range("D32").Value = '0.7234'
$ws.Range("E32").Value = '  -8.85%  '
$ws.Range("D33").Value = '2.860'
$ws.Range("E33").Value = '  -4.41%  '
$ws.Range("D34").Value = '4.318'
$ws.Range("E34").Value = '  -7.06%  '
$ws.Range("D35").Value = '1.084'
$ws.Range("E35").Value = '  -8.57%  '
$ws.Range("D36").Value = '1.001'
$ws.Range("E37").Value = '  -3.26%  '
$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D38").Value = '0.05130'
$ws.Range("E38").Value = '  -6.32%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").Value = '0.01899'
$ws.Range("E39").Value = '  -3.70%  '
$ws.Range("E40").Value = '  -4.40%  '
$ws.Range("D41").Value = '0.4959'
$ws.Range("E41").Value = '  -5.00%  '
$ws.Range("D42").Value = '2.597'
$ws.Range("E42").Value = '  -8.94%  '
$ws.Range("D43").Value = '8.038'
$ws.Range("E43").Value = '  -7.81%  '
$ws.Range("D44").Value = '5.910'
$ws.Range("E44").Value = '  -13.26%  '
$ws.Range("D45").Value = '105.00'
$ws.Range("E45").Value = '  -2.09%  '
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").Value = '10.22'
$ws.Range("E46").Value = '  -4.24%  '
$ws.Range("B47").Value = 'PaxDollar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D47").Value = '1.000'
$ws.Range("E47").Value = '  -0.08%  '
$ws.Range("D48").Value = '0.06293'
$ws.Range("E48").Value = '  -4.26%  '
$ws.Range("D49").Value = '0.4512'
$ws.Range("E49").Value = '  -6.60%  '
$ws.Range("D50").Value = '1.602'
$ws.Range("E50").Value = '  -4.67%  '
$ws.Range("D51").Value = '1.702'
$ws.Range("E51").Value = '  -8.85%  '
